$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 2351
$ws.Range("L3").Value = 2390
$ws.Range("E4").Value = 2052
$ws.Range("L4").Value = 647
$ws.Range("L6").Value = 2151
$ws.Range("E7").Value = 26057
$ws.Range("L7").Value = 7680

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L5").Value = 27
$ws.Range("E8").Value = 1913
$ws.Range("L8").Value = 490
$ws.Range("L11").Value = 135
$ws.Range("L15").Value = 54
$ws.Range("L16").Value = 15
$ws.Range("L18").Value = 54
$ws.Range("L19").Value = 220
$ws.Range("L20").Value = 196
$ws.Range("L25").Value = 44
$ws.Range("L27").Value = 79
$ws.Range("L29").Value = 392
$ws.Range("L33").Value = 348
$ws.Range("L37").Value = 280
$ws.Range("L42").Value = 238
$ws.Range("L44").Value = 56
$ws.Range("L47").Value = 60
$ws.Range("L48").Value = 108
$ws.Range("L51").Value = 88
$ws.Range("L52").Value = 153
$ws.Range("L53").Value = 97
$ws.Range("L57").Value = 33
$ws.Range("L60").Value = 45
$ws.Range("L66").Value = 17
$ws.Range("L67").Value = 280
$ws.Range("L68").Value = 22
$ws.Range("L72").Value = 32
$ws.Range("L76").Value = 88
$ws.Range("L79").Value = 210
$ws.Range("L80").Value = 28
$ws.Range("L83").Value = 188
$ws.Range("L85").Value = 398
$ws.Range("L87").Value = 23
$ws.Range("L88").Value = 105
$ws.Range("L89").Value = 96
$ws.Range("L91").Value = 111
$ws.Range("L94").Value = 90
$ws.Range("L96").Value = 77
$ws.Range("L99").Value = 122
$ws.Range("L100").Value = 12
$ws.Range("E101").Value = 26057
$ws.Range("L101").Value = 7680

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 19
$ws.Range("L4").Value = 11
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 77

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L3").Value = 42
$ws.Range("L7").Value = 135

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 96

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 163
$ws.Range("L4").Value = 32
$ws.Range("L7").Value = 398

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 49
$ws.Range("L7").Value = 153

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L3").Value = 165
$ws.Range("E4").Value = 124
$ws.Range("L4").Value = 37
$ws.Range("L6").Value = 130
$ws.Range("E7").Value = 1913
$ws.Range("L7").Value = 490

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L3").Value = 79
$ws.Range("L7").Value = 188

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L4").Value = 20
$ws.Range("L7").Value = 348

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 84
$ws.Range("L7").Value = 280

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 30
$ws.Range("L7").Value = 122

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 83
$ws.Range("L3").Value = 95
$ws.Range("L7").Value = 280

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 30
$ws.Range("L6").Value = 81

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L3").Value = 141
$ws.Range("L6").Value = 101
$ws.Range("L7").Value = 392

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L6").Value = 45
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 72
$ws.Range("L6").Value = 70
$ws.Range("L7").Value = 220

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 56

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 15
$ws.Range("L6").Value = 42
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 61
$ws.Range("L7").Value = 238

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 41
$ws.Range("L3").Value = 40
$ws.Range("L7").Value = 111

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 69
$ws.Range("L7").Value = 210

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 63
$ws.Range("L6").Value = 55
$ws.Range("L7").Value = 196

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L3").Value = 21
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("L6").Value = 5
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L6").Value = 14
$ws.Range("L7").Value = 60

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L6").Value = 12
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L6").Value = 6
$ws.Range("L7").Value = 17

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 105

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L3").Value = 8
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 27

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L6").Value = 21
$ws.Range("L7").Value = 79

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 25
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 22

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L2").Value = 12
$ws.Range("L7").Value = 33

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 45

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L6").Value = 7
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 28

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("L4").Value = 4
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 15

